$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap entire match rows (pairs got reordered in the source feed) ---
# Columns B and F:AC are swapped between each pair; A, C, D, E (id/div/date) stay put.
$ws.Range("B2").Value = 6983357
$ws.Range("B3").Value = 6979406
$ws.Range("F2").Value = "FK Cukaricki"
$ws.Range("F3").Value = "Spartak Subotica"
$ws.Range("G2").Value = "Radnicki Nis"
$ws.Range("G3").Value = "IMT Novi Belgrade"
$ws.Range("H2").Value = 2
$ws.Range("H3").Value = 2
$ws.Range("I2").Value = 0
$ws.Range("I3").Value = 1
$ws.Range("J2").Value = "H"
$ws.Range("J3").Value = "H"
$ws.Range("K2").Value = 1.571
$ws.Range("K3").Value = 1.909
$ws.Range("L2").Value = 3.6
$ws.Range("L3").Value = 3.2
$ws.Range("M2").Value = 5
$ws.Range("M3").Value = 3.6
$ws.Range("N2").Value = 1.533
$ws.Range("N3").Value = 1.909
$ws.Range("O2").Value = 3.6
$ws.Range("O3").Value = 3.1
$ws.Range("P2").Value = 5.25
$ws.Range("P3").Value = 3.8
$ws.Range("Q2").Value = -1
$ws.Range("Q3").Value = -0.5
$ws.Range("R2").Value = 2
$ws.Range("R3").Value = 1.975
$ws.Range("S2").Value = 1.8
$ws.Range("S3").Value = 1.825
$ws.Range("T2").Value = 2.5
$ws.Range("T3").Value = 2.25
$ws.Range("U2").Value = 1.925
$ws.Range("U3").Value = 1.9
$ws.Range("V2").Value = 1.875
$ws.Range("V3").Value = 1.9
$ws.Range("W2").Value = 0.5329999999999999
$ws.Range("W3").Value = 0.909
$ws.Range("X2").Value = -1
$ws.Range("X3").Value = -1
$ws.Range("Y2").Value = -1
$ws.Range("Y3").Value = -1
$ws.Range("Z2").Value = 1
$ws.Range("Z3").Value = 0.9750000000000001
$ws.Range("AA2").Value = -1
$ws.Range("AA3").Value = -1
$ws.Range("AB2").Value = -1
$ws.Range("AB3").Value = 0.8999999999999999
$ws.Range("AC2").Value = 0.875
$ws.Range("AC3").Value = -1
$ws.Range("B7").Value = 6978735
$ws.Range("B8").Value = 6979404
$ws.Range("F7").Value = "Red Star Belgrade"
$ws.Range("F8").Value = "FK Napredak"
$ws.Range("G7").Value = "Vojvodina"
$ws.Range("G8").Value = "FK Zeleznicar Pancevo"
$ws.Range("H7").Value = 5
$ws.Range("H8").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("I8").Value = 1
$ws.Range("J7").Value = "H"
$ws.Range("J8").Value = "A"
$ws.Range("K7").Value = 1.166
$ws.Range("K8").Value = 1.909
$ws.Range("L7").Value = 6
$ws.Range("L8").Value = 3.2
$ws.Range("M7").Value = 11
$ws.Range("M8").Value = 3.6
$ws.Range("N7").Value = 1.05
$ws.Range("N8").Value = 1.909
$ws.Range("O7").Value = 13
$ws.Range("O8").Value = 3.2
$ws.Range("P7").Value = 34
$ws.Range("P8").Value = 3.75
$ws.Range("Q7").Value = -3
$ws.Range("Q8").Value = -0.5
$ws.Range("R7").Value = 1.8
$ws.Range("R8").Value = 2
$ws.Range("S7").Value = 2
$ws.Range("S8").Value = 1.8
$ws.Range("T7").Value = 3.75
$ws.Range("T8").Value = 2.25
$ws.Range("U7").Value = 1.8
$ws.Range("U8").Value = 1.975
$ws.Range("V7").Value = 2
$ws.Range("V8").Value = 1.825
$ws.Range("W7").Value = 0.05000000000000004
$ws.Range("W8").Value = -1
$ws.Range("X7").Value = -1
$ws.Range("X8").Value = -1
$ws.Range("Y7").Value = -1
$ws.Range("Y8").Value = 2.75
$ws.Range("Z7").Value = 0.8
$ws.Range("Z8").Value = -1
$ws.Range("AA7").Value = -1
$ws.Range("AA8").Value = 0.8
$ws.Range("AB7").Value = 0.8
$ws.Range("AB8").Value = -1
$ws.Range("AC7").Value = -1
$ws.Range("AC8").Value = 0.825
$ws.Range("B19").Value = 7032917
$ws.Range("B20").Value = 7032914
$ws.Range("F19").Value = "FK Backa Topola"
$ws.Range("F20").Value = "FK Vozdovac"
$ws.Range("G19").Value = "FK Radnicki 1923"
$ws.Range("G20").Value = "FK Radnik Surdulica"
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("I19").Value = 0
$ws.Range("I20").Value = 1
$ws.Range("J19").Value = "H"
$ws.Range("J20").Value = "D"
$ws.Range("K19").Value = 1.5
$ws.Range("K20").Value = 2.2
$ws.Range("L19").Value = 3.75
$ws.Range("L20").Value = 3.1
$ws.Range("M19").Value = 6.5
$ws.Range("M20").Value = 3.2
$ws.Range("N19").Value = 1.444
$ws.Range("N20").Value = 2.05
$ws.Range("O19").Value = 4
$ws.Range("O20").Value = 3.1
$ws.Range("P19").Value = 6.5
$ws.Range("P20").Value = 3.5
$ws.Range("Q19").Value = -1.25
$ws.Range("Q20").Value = -0.25
$ws.Range("R19").Value = 1.975
$ws.Range("R20").Value = 1.75
$ws.Range("S19").Value = 1.825
$ws.Range("S20").Value = 2.05
$ws.Range("T19").Value = 2.75
$ws.Range("T20").Value = 2
$ws.Range("U19").Value = 1.95
$ws.Range("U20").Value = 1.775
$ws.Range("V19").Value = 1.85
$ws.Range("V20").Value = 2.025
$ws.Range("W19").Value = 0.444
$ws.Range("W20").Value = -1
$ws.Range("X19").Value = -1
$ws.Range("X20").Value = 2.1
$ws.Range("Y19").Value = -1
$ws.Range("Y20").Value = -1
$ws.Range("Z19").Value = -0.5
$ws.Range("Z20").Value = -0.5
$ws.Range("AA19").Value = 0.4125
$ws.Range("AA20").Value = 0.5249999999999999
$ws.Range("AB19").Value = -1
$ws.Range("AB20").Value = 0
$ws.Range("AC19").Value = 0.8500000000000001
$ws.Range("AC20").Value = 0
$ws.Range("B43").Value = 6978740
$ws.Range("B44").Value = 6979440
$ws.Range("F43").Value = "Red Star Belgrade"
$ws.Range("F44").Value = "Javor Ivanjica"
$ws.Range("G43").Value = "FK Novi Pazar"
$ws.Range("G44").Value = "Radnicki Nis"
$ws.Range("H43").Value = 2
$ws.Range("H44").Value = 1
$ws.Range("I43").Value = 1
$ws.Range("I44").Value = 0
$ws.Range("J43").Value = "H"
$ws.Range("J44").Value = "H"
$ws.Range("K43").Value = 1.062
$ws.Range("K44").Value = 2.3
$ws.Range("L43").Value = 13
$ws.Range("L44").Value = 3.2
$ws.Range("M43").Value = 23
$ws.Range("M44").Value = 2.875
$ws.Range("N43").Value = 1.025
$ws.Range("N44").Value = 2.5
$ws.Range("O43").Value = 19
$ws.Range("O44").Value = 3.25
$ws.Range("P43").Value = 41
$ws.Range("P44").Value = 2.6
$ws.Range("Q43").Value = -3.75
$ws.Range("Q44").Value = 0
$ws.Range("R43").Value = 1.825
$ws.Range("R44").Value = 1.85
$ws.Range("S43").Value = 1.975
$ws.Range("S44").Value = 1.95
$ws.Range("T43").Value = 4.5
$ws.Range("T44").Value = 2.25
$ws.Range("U43").Value = 1.975
$ws.Range("U44").Value = 1.9
$ws.Range("V43").Value = 1.825
$ws.Range("V44").Value = 1.9
$ws.Range("W43").Value = 0.02499999999999991
$ws.Range("W44").Value = 1.5
$ws.Range("X43").Value = -1
$ws.Range("X44").Value = -1
$ws.Range("Y43").Value = -1
$ws.Range("Y44").Value = -1
$ws.Range("Z43").Value = -1
$ws.Range("Z44").Value = 0.8500000000000001
$ws.Range("AA43").Value = 0.9750000000000001
$ws.Range("AA44").Value = -1
$ws.Range("AB43").Value = -1
$ws.Range("AB44").Value = -1
$ws.Range("AC43").Value = 0.825
$ws.Range("AC44").Value = 0.8999999999999999
$ws.Range("B50").Value = 6979447
$ws.Range("B51").Value = 6979449
$ws.Range("F50").Value = "FK Vozdovac"
$ws.Range("F51").Value = "FK Backa Topola"
$ws.Range("G50").Value = "Mladost Lucani"
$ws.Range("G51").Value = "FK Zeleznicar Pancevo"
$ws.Range("H50").Value = 3
$ws.Range("H51").Value = 6
$ws.Range("I50").Value = 1
$ws.Range("I51").Value = 3
$ws.Range("J50").Value = "H"
$ws.Range("J51").Value = "H"
$ws.Range("K50").Value = 1.909
$ws.Range("K51").Value = 1.25
$ws.Range("L50").Value = 3.25
$ws.Range("L51").Value = 5
$ws.Range("M50").Value = 3.5
$ws.Range("M51").Value = 9
$ws.Range("N50").Value = 1.95
$ws.Range("N51").Value = 1.285
$ws.Range("O50").Value = 3.1
$ws.Range("O51").Value = 4.5
$ws.Range("P50").Value = 3.6
$ws.Range("P51").Value = 9.5
$ws.Range("Q50").Value = -0.5
$ws.Range("Q51").Value = -1.75
$ws.Range("R50").Value = 2
$ws.Range("R51").Value = 2
$ws.Range("S50").Value = 1.8
$ws.Range("S51").Value = 1.8
$ws.Range("T50").Value = 2.25
$ws.Range("T51").Value = 2.75
$ws.Range("U50").Value = 1.9
$ws.Range("U51").Value = 1.875
$ws.Range("V50").Value = 1.9
$ws.Range("V51").Value = 1.925
$ws.Range("W50").Value = 0.95
$ws.Range("W51").Value = 0.2849999999999999
$ws.Range("X50").Value = -1
$ws.Range("X51").Value = -1
$ws.Range("Y50").Value = -1
$ws.Range("Y51").Value = -1
$ws.Range("Z50").Value = 1
$ws.Range("Z51").Value = 1
$ws.Range("AA50").Value = -1
$ws.Range("AA51").Value = -1
$ws.Range("AB50").Value = 0.8999999999999999
$ws.Range("AB51").Value = 0.875
$ws.Range("AC50").Value = -1
$ws.Range("AC51").Value = -1
$ws.Range("B61").Value = 6979453
$ws.Range("B62").Value = 6979458
$ws.Range("F61").Value = "Vojvodina"
$ws.Range("F62").Value = "Mladost Lucani"
$ws.Range("G61").Value = "Javor Ivanjica"
$ws.Range("G62").Value = "FK Novi Pazar"
$ws.Range("H61").Value = 2
$ws.Range("H62").Value = 2
$ws.Range("I61").Value = 1
$ws.Range("I62").Value = 0
$ws.Range("J61").Value = "H"
$ws.Range("J62").Value = "H"
$ws.Range("K61").Value = 1.4
$ws.Range("K62").Value = 2.3
$ws.Range("L61").Value = 4.1
$ws.Range("L62").Value = 3.2
$ws.Range("M61").Value = 7
$ws.Range("M62").Value = 2.875
$ws.Range("N61").Value = 1.363
$ws.Range("N62").Value = 2.375
$ws.Range("O61").Value = 4.5
$ws.Range("O62").Value = 3
$ws.Range("P61").Value = 7
$ws.Range("P62").Value = 2.9
$ws.Range("Q61").Value = -1.25
$ws.Range("Q62").Value = -0.25
$ws.Range("R61").Value = 1.825
$ws.Range("R62").Value = 2.025
$ws.Range("S61").Value = 1.975
$ws.Range("S62").Value = 1.775
$ws.Range("T61").Value = 2.75
$ws.Range("T62").Value = 2
$ws.Range("U61").Value = 1.975
$ws.Range("U62").Value = 1.75
$ws.Range("V61").Value = 1.825
$ws.Range("V62").Value = 2.05
$ws.Range("W61").Value = 0.363
$ws.Range("W62").Value = 1.375
$ws.Range("X61").Value = -1
$ws.Range("X62").Value = -1
$ws.Range("Y61").Value = -1
$ws.Range("Y62").Value = -1
$ws.Range("Z61").Value = -0.5
$ws.Range("Z62").Value = 1.025
$ws.Range("AA61").Value = 0.4875
$ws.Range("AA62").Value = -1
$ws.Range("AB61").Value = 0.4875
$ws.Range("AB62").Value = 0
$ws.Range("AC61").Value = -0.5
$ws.Range("AC62").Value = 0
$ws.Range("B84").Value = 6979481
$ws.Range("B85").Value = 6979484
$ws.Range("F84").Value = "Vojvodina"
$ws.Range("F85").Value = "Mladost Lucani"
$ws.Range("G84").Value = "FK Radnik Surdulica"
$ws.Range("G85").Value = "Radnicki Nis"
$ws.Range("H84").Value = 3
$ws.Range("H85").Value = 1
$ws.Range("I84").Value = 0
$ws.Range("I85").Value = 2
$ws.Range("J84").Value = "H"
$ws.Range("J85").Value = "A"
$ws.Range("K84").Value = 2.25
$ws.Range("K85").Value = 2
$ws.Range("L84").Value = 3
$ws.Range("L85").Value = 3.25
$ws.Range("M84").Value = 3
$ws.Range("M85").Value = 3.25
$ws.Range("N84").Value = 1.363
$ws.Range("N85").Value = 2.55
$ws.Range("O84").Value = 3.8
$ws.Range("O85").Value = 3.2
$ws.Range("P84").Value = 9
$ws.Range("P85").Value = 2.55
$ws.Range("Q84").Value = -1.25
$ws.Range("Q85").Value = 0
$ws.Range("R84").Value = 1.825
$ws.Range("R85").Value = 1.9
$ws.Range("S84").Value = 1.975
$ws.Range("S85").Value = 1.9
$ws.Range("T84").Value = 2.5
$ws.Range("T85").Value = 2.25
$ws.Range("U84").Value = 1.9
$ws.Range("U85").Value = 1.875
$ws.Range("V84").Value = 1.9
$ws.Range("V85").Value = 1.925
$ws.Range("W84").Value = 0.363
$ws.Range("W85").Value = -1
$ws.Range("X84").Value = -1
$ws.Range("X85").Value = -1
$ws.Range("Y84").Value = -1
$ws.Range("Y85").Value = 1.55
$ws.Range("Z84").Value = 0.825
$ws.Range("Z85").Value = -1
$ws.Range("AA84").Value = -1
$ws.Range("AA85").Value = 0.8999999999999999
$ws.Range("AB84").Value = 0.8999999999999999
$ws.Range("AB85").Value = 0.875
$ws.Range("AC84").Value = -1
$ws.Range("AC85").Value = -1
$ws.Range("B90").Value = 6979491
$ws.Range("B91").Value = 6978747
$ws.Range("F90").Value = "Radnicki Nis"
$ws.Range("F91").Value = "IMT Novi Belgrade"
$ws.Range("G90").Value = "Spartak Subotica"
$ws.Range("G91").Value = "Red Star Belgrade"
$ws.Range("H90").Value = 1
$ws.Range("H91").Value = 1
$ws.Range("I90").Value = 1
$ws.Range("I91").Value = 2
$ws.Range("J90").Value = "D"
$ws.Range("J91").Value = "A"
$ws.Range("K90").Value = 1.95
$ws.Range("K91").Value = 8
$ws.Range("L90").Value = 3.25
$ws.Range("L91").Value = 5.25
$ws.Range("M90").Value = 3.7
$ws.Range("M91").Value = 1.285
$ws.Range("N90").Value = 1.65
$ws.Range("N91").Value = 15
$ws.Range("O90").Value = 3.5
$ws.Range("O91").Value = 7.5
$ws.Range("P90").Value = 5
$ws.Range("P91").Value = 1.125
$ws.Range("Q90").Value = -0.75
$ws.Range("Q91").Value = 2.25
$ws.Range("R90").Value = 1.825
$ws.Range("R91").Value = 1.975
$ws.Range("S90").Value = 1.975
$ws.Range("S91").Value = 1.825
$ws.Range("T90").Value = 2.5
$ws.Range("T91").Value = 3.5
$ws.Range("U90").Value = 2
$ws.Range("U91").Value = 1.825
$ws.Range("V90").Value = 1.8
$ws.Range("V91").Value = 1.975
$ws.Range("W90").Value = -1
$ws.Range("W91").Value = -1
$ws.Range("X90").Value = 2.5
$ws.Range("X91").Value = -1
$ws.Range("Y90").Value = -1
$ws.Range("Y91").Value = 0.125
$ws.Range("Z90").Value = -1
$ws.Range("Z91").Value = 0.9750000000000001
$ws.Range("AA90").Value = 0.9750000000000001
$ws.Range("AA91").Value = -1
$ws.Range("AB90").Value = -1
$ws.Range("AB91").Value = -1
$ws.Range("AC90").Value = 0.8
$ws.Range("AC91").Value = 0.9750000000000001
$ws.Range("B124").Value = 6979516
$ws.Range("B125").Value = 6979522
$ws.Range("F124").Value = "Partizan Belgrade"
$ws.Range("F125").Value = "Mladost Lucani"
$ws.Range("G124").Value = "Vojvodina"
$ws.Range("G125").Value = "FK Zeleznicar Pancevo"
$ws.Range("H124").Value = 3
$ws.Range("H125").Value = 1
$ws.Range("I124").Value = 1
$ws.Range("I125").Value = 0
$ws.Range("J124").Value = "H"
$ws.Range("J125").Value = "H"
$ws.Range("K124").Value = 1.5
$ws.Range("K125").Value = 2.15
$ws.Range("L124").Value = 4
$ws.Range("L125").Value = 3.25
$ws.Range("M124").Value = 5.5
$ws.Range("M125").Value = 3.1
$ws.Range("N124").Value = 1.444
$ws.Range("N125").Value = 2.1
$ws.Range("O124").Value = 4.2
$ws.Range("O125").Value = 3.4
$ws.Range("P124").Value = 6
$ws.Range("P125").Value = 3.1
$ws.Range("Q124").Value = -1.25
$ws.Range("Q125").Value = -0.25
$ws.Range("R124").Value = 2.025
$ws.Range("R125").Value = 1.85
$ws.Range("S124").Value = 1.775
$ws.Range("S125").Value = 1.95
$ws.Range("T124").Value = 2.75
$ws.Range("T125").Value = 2.25
$ws.Range("U124").Value = 1.775
$ws.Range("U125").Value = 1.775
$ws.Range("V124").Value = 2.025
$ws.Range("V125").Value = 2.025
$ws.Range("W124").Value = 0.444
$ws.Range("W125").Value = 1.1
$ws.Range("X124").Value = -1
$ws.Range("X125").Value = -1
$ws.Range("Y124").Value = -1
$ws.Range("Y125").Value = -1
$ws.Range("Z124").Value = 1.025
$ws.Range("Z125").Value = 0.8500000000000001
$ws.Range("AA124").Value = -1
$ws.Range("AA125").Value = -1
$ws.Range("AB124").Value = 0.7749999999999999
$ws.Range("AB125").Value = -1
$ws.Range("AC124").Value = -1
$ws.Range("AC125").Value = 1.025
$ws.Range("B190").Value = 6979566
$ws.Range("B191").Value = 7921659
$ws.Range("F190").Value = "Mladost Lucani"
$ws.Range("F191").Value = "Vojvodina"
$ws.Range("G190").Value = "FK Cukaricki"
$ws.Range("G191").Value = "FK Vozdovac"
$ws.Range("H190").Value = 1
$ws.Range("H191").Value = 2
$ws.Range("I190").Value = 0
$ws.Range("I191").Value = 1
$ws.Range("J190").Value = "H"
$ws.Range("J191").Value = "H"
$ws.Range("K190").Value = 3.4
$ws.Range("K191").Value = 1.6
$ws.Range("L190").Value = 3.3
$ws.Range("L191").Value = 3.6
$ws.Range("M190").Value = 1.95
$ws.Range("M191").Value = 4.75
$ws.Range("N190").Value = 2.7
$ws.Range("N191").Value = 1.45
$ws.Range("O190").Value = 3
$ws.Range("O191").Value = 3.8
$ws.Range("P190").Value = 2.45
$ws.Range("P191").Value = 6
$ws.Range("Q190").Value = 0
$ws.Range("Q191").Value = -1
$ws.Range("R190").Value = 2
$ws.Range("R191").Value = 1.8
$ws.Range("S190").Value = 1.8
$ws.Range("S191").Value = 2
$ws.Range("T190").Value = 2
$ws.Range("T191").Value = 2.5
$ws.Range("U190").Value = 1.75
$ws.Range("U191").Value = 1.85
$ws.Range("V190").Value = 2.05
$ws.Range("V191").Value = 1.95
$ws.Range("W190").Value = 1.7
$ws.Range("W191").Value = 0.45
$ws.Range("X190").Value = -1
$ws.Range("X191").Value = -1
$ws.Range("Y190").Value = -1
$ws.Range("Y191").Value = -1
$ws.Range("Z190").Value = 1
$ws.Range("Z191").Value = 0
$ws.Range("AA190").Value = -1
$ws.Range("AA191").Value = 0
$ws.Range("AB190").Value = -1
$ws.Range("AB191").Value = 0.8500000000000001
$ws.Range("AC190").Value = 1.05
$ws.Range("AC191").Value = -1

# --- Rewrite rows 234-236 with updated match data/results ---
$ws.Range("B234").Value = 6998838
$ws.Range("E234").Value = 45396.47916666666
$ws.Range("F234").Value = "Vojvodina"
$ws.Range("G234").Value = "FK Backa Topola"
$ws.Range("H234").Value = 3
$ws.Range("I234").Value = 2
$ws.Range("J234").Value = "H"
$ws.Range("K234").Value = 2.3
$ws.Range("L234").Value = 3.4
$ws.Range("M234").Value = 2.75
$ws.Range("N234").Value = 2.375
$ws.Range("O234").Value = 3.4
$ws.Range("P234").Value = 2.7
$ws.Range("Q234").Value = 0
$ws.Range("R234").Value = 1.775
$ws.Range("S234").Value = 2.025
$ws.Range("T234").Value = 2.5
$ws.Range("U234").Value = 1.85
$ws.Range("V234").Value = 1.95
$ws.Range("W234").Value = 1.375
$ws.Range("X234").Value = -1
$ws.Range("Y234").Value = -1
$ws.Range("Z234").Value = 0.7749999999999999
$ws.Range("AA234").Value = -1
$ws.Range("AB234").Value = 0.8500000000000001
$ws.Range("AC234").Value = -1
$ws.Range("B235").Value = 6979612
$ws.Range("E235").Value = 45396.47916666666
$ws.Range("F235").Value = "FK Napredak"
$ws.Range("G235").Value = "FK Radnik Surdulica"
$ws.Range("H235").Value = 0
$ws.Range("I235").Value = 2
$ws.Range("J235").Value = "A"
$ws.Range("K235").Value = 1.5
$ws.Range("L235").Value = 4.1
$ws.Range("M235").Value = 5.5
$ws.Range("N235").Value = 1.833
$ws.Range("O235").Value = 3.4
$ws.Range("P235").Value = 4
$ws.Range("Q235").Value = -0.5
$ws.Range("R235").Value = 1.825
$ws.Range("S235").Value = 1.975
$ws.Range("T235").Value = 2
$ws.Range("U235").Value = 1.875
$ws.Range("V235").Value = 1.925
$ws.Range("W235").Value = -1
$ws.Range("X235").Value = -1
$ws.Range("Y235").Value = 3
$ws.Range("Z235").Value = -1
$ws.Range("AA235").Value = 0.9750000000000001
$ws.Range("AB235").Value = 0
$ws.Range("AC235").Value = 0
$ws.Range("B236").Value = 6979613
$ws.Range("E236").Value = 45396.5
$ws.Range("F236").Value = "Crvena Zvezda"
$ws.Range("G236").Value = "FK Zeleznicar Pancevo"
$ws.Range("H236").Value = 3
$ws.Range("I236").Value = 0
$ws.Range("J236").Value = "H"
$ws.Range("K236").Value = 1.1
$ws.Range("L236").Value = 8.5
$ws.Range("M236").Value = 15
$ws.Range("N236").Value = 1.083
$ws.Range("O236").Value = 10
$ws.Range("P236").Value = 21
$ws.Range("Q236").Value = -2.75
$ws.Range("R236").Value = 1.85
$ws.Range("S236").Value = 1.95
$ws.Range("T236").Value = 4
$ws.Range("U236").Value = 1.975
$ws.Range("V236").Value = 1.825
$ws.Range("W236").Value = 0.08299999999999996
$ws.Range("X236").Value = -1
$ws.Range("Y236").Value = -1
$ws.Range("Z236").Value = 0.425
$ws.Range("AA236").Value = -0.5
$ws.Range("AB236").Value = -1
$ws.Range("AC236").Value = 0.825
